$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 3330.5
$ws.Range("I29").Value = 3330.5
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 9991.5
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -9710.5
$ws.Range("N29").ClearContents()
$ws.Range("H38").Value = 2403.9092
$ws.Range("J38").Value = 3502.25
$ws.Range("L38").Value = 10506.75
$ws.Range("N38").Value = -11250.75
$ws.Range("H39").Value = 891.5714
$ws.Range("I39").Value = 921.6923
$ws.Range("J39").Value = 500
$ws.Range("K39").Value = 2765.0769
$ws.Range("L39").Value = 1500
$ws.Range("M39").Value = -2469.0769
$ws.Range("N39").Value = -2092
$ws.Range("H43").Value = 460
$ws.Range("I43").Value = 400
$ws.Range("J43").Value = 490
$ws.Range("K43").Value = 400
$ws.Range("L43").Value = 490
$ws.Range("M43").Value = -331
$ws.Range("N43").Value = -628
$ws.Range("H137").Value = 1499.3214
$ws.Range("I137").Value = 944.55554
$ws.Range("J137").Value = 1762.1052
$ws.Range("K137").Value = 2833.66662
$ws.Range("L137").Value = 5286.3156
$ws.Range("M137").Value = -283.66662
$ws.Range("N137").Value = -10386.3156
$ws.Range("H138").Value = 1823.71
$ws.Range("J138").Value = 2373.5735
$ws.Range("L138").Value = 7120.720499999999
$ws.Range("N138").Value = -17400.7205

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5922.64
$ws.Range("I2").Value = 7277.222
$ws.Range("J2").Value = 2439.4285
$ws.Range("K2").Value = 7277.222
$ws.Range("L2").Value = 2439.4285
$ws.Range("M2").Value = -7164.222
$ws.Range("N2").Value = -2665.4285
$ws.Range("H32").Value = 2555.4624
$ws.Range("I32").Value = 1735.0469
$ws.Range("J32").Value = 5837.125
$ws.Range("K32").Value = 1735.0469
$ws.Range("L32").Value = 5837.125
$ws.Range("M32").Value = -1448.0469
$ws.Range("N32").Value = -6411.125
$ws.Range("H61").Value = 3344.1052
$ws.Range("I61").Value = 1844.5714
$ws.Range("J61").Value = 4218.8335
$ws.Range("K61").Value = 1844.5714
$ws.Range("L61").Value = 4218.8335
$ws.Range("M61").Value = -1632.5714
$ws.Range("N61").Value = -4642.8335
$ws.Range("H74").Value = 1369.1875
$ws.Range("I74").Value = 1416.6666
$ws.Range("J74").Value = 1278.5454
$ws.Range("K74").Value = 1416.6666
$ws.Range("L74").Value = 1278.5454
$ws.Range("M74").Value = -542.6666
$ws.Range("N74").Value = -3026.5454
$ws.Range("H77").Value = 1369.1875
$ws.Range("I77").Value = 1416.6666
$ws.Range("J77").Value = 1278.5454
$ws.Range("K77").Value = 7083.333000000001
$ws.Range("L77").Value = 6392.727
$ws.Range("M77").Value = -2715.333000000001
$ws.Range("N77").Value = -15128.727
$ws.Range("H116").Value = 5922.64
$ws.Range("I116").Value = 7277.222
$ws.Range("J116").Value = 2439.4285
$ws.Range("K116").Value = 7277.222
$ws.Range("L116").Value = 2439.4285
$ws.Range("M116").Value = -4983.222
$ws.Range("N116").Value = -7027.4285
$ws.Range("H122").Value = 2104.4
$ws.Range("I122").Value = 1680
$ws.Range("J122").Value = 2528.8
$ws.Range("K122").Value = 5040
$ws.Range("L122").Value = 7586.400000000001
$ws.Range("M122").Value = -2590
$ws.Range("N122").Value = -12486.4
$ws.Range("H136").Value = 3344.1052
$ws.Range("I136").Value = 1844.5714
$ws.Range("J136").Value = 4218.8335
$ws.Range("K136").Value = 5533.7142
$ws.Range("L136").Value = 12656.5005
$ws.Range("M136").Value = -2983.7142
$ws.Range("N136").Value = -17756.5005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5922.64
$ws.Range("I3").Value = 7277.222
$ws.Range("J3").Value = 2439.4285
$ws.Range("K3").Value = 7277.222
$ws.Range("L3").Value = 2439.4285
$ws.Range("M3").Value = -7163.222
$ws.Range("N3").Value = -2667.4285
$ws.Range("H134").Value = 2676.3076
$ws.Range("I134").Value = 1777.4667
$ws.Range("J134").Value = 5672.4443
$ws.Range("K134").Value = 5332.4001
$ws.Range("L134").Value = 17017.3329
$ws.Range("M134").Value = -2797.4001
$ws.Range("N134").Value = -22087.3329

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1487.6234
$ws.Range("I31").Value = 1126.625
$ws.Range("J31").Value = 1744.3334
$ws.Range("K31").Value = 1126.625
$ws.Range("L31").Value = 1744.3334
$ws.Range("M31").Value = -831.625
$ws.Range("N31").Value = -2334.3334
$ws.Range("H34").Value = 1487.6234
$ws.Range("I34").Value = 1126.625
$ws.Range("J34").Value = 1744.3334
$ws.Range("K34").Value = 1126.625
$ws.Range("L34").Value = 1744.3334
$ws.Range("M34").Value = -924.625
$ws.Range("N34").Value = -2148.3334
$ws.Range("H58").Value = 2261.72
$ws.Range("I58").Value = 876.2
$ws.Range("J58").Value = 3185.4
$ws.Range("K58").Value = 876.2
$ws.Range("L58").Value = 3185.4
$ws.Range("M58").Value = -673.2
$ws.Range("N58").Value = -3591.4
$ws.Range("H105").Value = 487.93332
$ws.Range("I105").Value = 430.6
$ws.Range("J105").Value = 602.6
$ws.Range("K105").Value = 430.6
$ws.Range("L105").Value = 602.6
$ws.Range("M105").Value = 1316.4
$ws.Range("N105").Value = -4096.6
$ws.Range("H132").Value = 2820.7188
$ws.Range("I132").Value = 2086.75
$ws.Range("J132").Value = 4044
$ws.Range("K132").Value = 6260.25
$ws.Range("L132").Value = 12132
$ws.Range("M132").Value = -3730.25
$ws.Range("N132").Value = -17192
$ws.Range("H134").Value = 2784.6155
$ws.Range("I134").Value = 1223
$ws.Range("K134").Value = 3669
$ws.Range("M134").Value = -1134
$ws.Range("H136").Value = 2261.72
$ws.Range("I136").Value = 876.2
$ws.Range("J136").Value = 3185.4
$ws.Range("K136").Value = 2628.6
$ws.Range("L136").Value = 9556.200000000001
$ws.Range("M136").Value = -78.60000000000036
$ws.Range("N136").Value = -14656.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 22800
$ws.Range("J33").Value = 22800
$ws.Range("L33").Value = 22800
$ws.Range("N33").Value = -23304

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 29257.75
$ws.Range("J42").Value = 8515.5
$ws.Range("L42").Value = 8515.5
$ws.Range("N42").Value = -9641.5
$ws.Range("H49").Value = 29257.75
$ws.Range("J49").Value = 8515.5
$ws.Range("L49").Value = 8515.5
$ws.Range("N49").Value = -8809.5
$ws.Range("H55").Value = 358.46155
$ws.Range("I55").Value = 309.875
$ws.Range("J55").Value = 436.2
$ws.Range("K55").Value = 309.875
$ws.Range("L55").Value = 436.2
$ws.Range("M55").Value = -136.875
$ws.Range("N55").Value = -782.2
$ws.Range("H132").Value = 2580.75
$ws.Range("I132").Value = 1795.5758
$ws.Range("J132").Value = 4308.1333
$ws.Range("K132").Value = 5386.7274
$ws.Range("L132").Value = 12924.3999
$ws.Range("M132").Value = -2856.7274
$ws.Range("N132").Value = -17984.3999
$ws.Range("H136").Value = 3448.4443
$ws.Range("I136").Value = 1350.0952
$ws.Range("J136").Value = 6386.1333
$ws.Range("K136").Value = 4050.2856
$ws.Range("L136").Value = 19158.3999
$ws.Range("M136").Value = -1500.2856
$ws.Range("N136").Value = -24258.3999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 40020000
$ws.Range("J2").Value = 47997.668
$ws.Range("L2").Value = 47997.668
$ws.Range("N2").Value = -48221.668
$ws.Range("H122").Value = 78145.84
$ws.Range("I122").Value = 143814.42
$ws.Range("J122").Value = 1532.5
$ws.Range("K122").Value = 431443.26
$ws.Range("L122").Value = 4597.5
$ws.Range("M122").Value = -428993.26
$ws.Range("N122").Value = -9497.5
$ws.Range("H132").Value = 19233134
$ws.Range("I132").Value = 23811412
$ws.Range("K132").Value = 71434236
$ws.Range("M132").Value = -71431706
$ws.Range("H136").Value = 12860332
$ws.Range("I136").Value = 16717319
$ws.Range("J136").Value = 3708.8333
$ws.Range("K136").Value = 50151957
$ws.Range("L136").Value = 11126.4999
$ws.Range("M136").Value = -50149407
$ws.Range("N136").Value = -16226.4999
